# cryptos.xlsx refresh: updates Price (col D) and Volume(1h) (col E) values
# scraped from coinranking.com, plus rows 49/50 which were re-sorted so the
# Coin/Link/Price/Volume columns swap between EnergySwap and Algorand.
#
# Column D cells are free-text price strings (e.g. '26.783.82', '1.00',
# '0.0845') rather than numbers, matching the source sheet's inlineStr cells.
# Assigning a bare numeric-looking string to Range.Value lets Excel's COM
# layer auto-coerce it to a real number (dropping formatting like trailing
# zeros / leading zeros), so those assignments use a leading apostrophe
# (the standard Excel 'force text' prefix) to keep them as literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '26.783.82'
$ws.Range('E2').Value = '  -0.20%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.637.71'
$ws.Range('E3').Value = '  -0.56%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  -0.55%  '

# Row 5: BNB
$ws.Range('D5').Value = '''219.06'
$ws.Range('E5').Value = '  +0.73%  '

# Row 6: XRP
$ws.Range('E6').Value = '  -0.90%  '

# Row 7: USDC
$ws.Range('D7').Value = '''1.01'
$ws.Range('E7').Value = '  -0.55%  '

# Row 8: Cardano
$ws.Range('E8').Value = '  -0.58%  '

# Row 9: Dogecoin
$ws.Range('E9').Value = '  -0.99%  '

# Row 10: Solana
$ws.Range('D10').Value = '''19.18'
$ws.Range('E10').Value = '  -0.16%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.0845'
$ws.Range('E11').Value = '  +0.41%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.866.92'
$ws.Range('E12').Value = '  -0.48%  '

# Row 13: WrappedEther
$ws.Range('D13').Value = '1.644.42'
$ws.Range('E13').Value = '  -0.22%  '

# Row 14: Polkadot
$ws.Range('E14').Value = '  -1.35%  '

# Row 15: Polygon
$ws.Range('E15').Value = '  -0.72%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '26.787.10'

# Row 18: ShibaInu
$ws.Range('E18').Value = '  -0.89%  '

# Row 19: BitcoinCash
$ws.Range('D19').Value = '''214.89'
$ws.Range('E19').Value = '  +0.09%  '

# Row 20: Dai
$ws.Range('D20').Value = '''1.01'
$ws.Range('E20').Value = '  -0.64%  '

# Row 21: Uniswap
$ws.Range('E21').Value = '  -0.11%  '

# Row 22: Chainlink
$ws.Range('E22').Value = '  -0.19%  '

# Row 23: Toncoin
$ws.Range('E23').Value = '  -3.11%  '

# Row 24: Avalanche
$ws.Range('D24').Value = '''9.10'
$ws.Range('E24').Value = '  -2.83%  '

# Row 25: Monero
$ws.Range('E25').Value = '  +1.70%  '

# Row 26: BinanceUSD
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  -0.94%  '

# Row 27: Stellar
$ws.Range('E27').Value = '  -0.18%  '

# Row 28: Cosmos
$ws.Range('D28').Value = '''7.04'
$ws.Range('E28').Value = '  -0.64%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '''15.68'
$ws.Range('E29').Value = '  -0.16%  '

# Row 30: Hedera
$ws.Range('D30').Value = '''0.0506'
$ws.Range('E30').Value = '  -1.61%  '

# Row 31: PancakeSwap
$ws.Range('E31').Value = '  +1.48%  '

# Row 32: Filecoin
$ws.Range('D32').Value = '''3.38'
$ws.Range('E32').Value = '  +2.08%  '

# Row 33: InternetComputer(DFINITY)
$ws.Range('E33').Value = '  -0.04%  '

# Row 34: LidoDAOToken
$ws.Range('E34').Value = '  +0.03%  '

# Row 35: Maker
$ws.Range('D35').Value = '1.262.02'
$ws.Range('E35').Value = '  -1.94%  '

# Row 36: HuobiToken
$ws.Range('E36').Value = '  +0.00%  '

# Row 37: VeChain
$ws.Range('E37').Value = '  -0.03%  '

# Row 38: ImmutableX
$ws.Range('D38').Value = '''0.527'
$ws.Range('E38').Value = '  -1.94%  '

# Row 39: ARBITRUM
$ws.Range('E39').Value = '  -1.94%  '

# Row 40: PaxDollar
$ws.Range('E40').Value = '  -0.48%  '

# Row 41: TrustWalletToken
$ws.Range('D41').Value = '''0.805'
$ws.Range('E41').Value = '  -1.04%  '

# Row 42: FraxShare
$ws.Range('E42').Value = '  -0.53%  '

# Row 43: RocketPoolETH
$ws.Range('D43').Value = '1.776.97'
$ws.Range('E43').Value = '  -1.08%  '

# Row 44: MXToken
$ws.Range('E44').Value = '  -4.80%  '

# Row 45: Quant
$ws.Range('D45').Value = '''92.09'
$ws.Range('E45').Value = '  +0.61%  '

# Row 46: Aave
$ws.Range('D46').Value = '''59.94'
$ws.Range('E46').Value = '  +0.56%  '

# Row 47: RenderToken
$ws.Range('D47').Value = '''1.58'
$ws.Range('E47').Value = '  -1.23%  '

# Row 48: Cronos
$ws.Range('E48').Value = '  -0.72%  '

# Row 49: EnergySwap->Algorand
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '''0.0960'
$ws.Range('E49').Value = '  -1.34%  '

# Row 50: Algorand->EnergySwap
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = '''7.50'
$ws.Range('E50').Value = '  -2.00%  '

# Row 51: USDD
$ws.Range('D51').Value = '''1.01'
$ws.Range('E51').Value = '  -0.61%  '
